$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append a new observation record as row 3 (same species/site group as the
# existing row 2 "Spillkråka" / Dryocopus martius record, but a new find
# with its own id, nest note, coordinates and reporter).

$ws.Range("A3").Value = 131242840
$ws.Range("B3").Value = 57881
$ws.Range("D3").Value = "NT"
$ws.Range("E3").Value = 100049
$ws.Range("F3").Value = "Spillkråka"
$ws.Range("G3").Value = "Dryocopus martius"
$ws.Range("H3").Value = "(Linnaeus, 1758)"

# A few columns have no text for this record but still carry a (blank)
# cell in the source data - touching the number format materializes an
# entry for the cell even though its value stays empty.
$ws.Cells.Item(3, 9).NumberFormat = "General"   # I3  (Antal)
$ws.Cells.Item(3, 11).NumberFormat = "General"  # K3  (Ålder-Stadium)
$ws.Cells.Item(3, 12).NumberFormat = "General"  # L3  (Kön)

$ws.Range("M3").Value = "gammalt bo"

$ws.Cells.Item(3, 14).NumberFormat = "General"  # N3  (Metod)

$ws.Range("P3").Value = "Lille-Väktor, Boh"
$ws.Range("Q3").Value = 327428
$ws.Range("R3").Value = 6453550
$ws.Range("S3").Value = 5
$ws.Range("T3").Value = "Västra Götaland"
$ws.Range("U3").Value = "Lilla Edet"
$ws.Range("V3").Value = "Bohuslän"
$ws.Range("W3").Value = "Hjärtum"

# Startdatum / Slutdatum are stored as plain text (not Excel date serials)
# in this sheet, so force text with a leading apostrophe, matching how
# this date-looking text was entered originally.
$ws.Range("Y3").Value = "'2026-01-31"
$ws.Range("AA3").Value = "'2026-01-31"

$ws.Range("AD3").Value = $false
$ws.Range("AE3").Value = $false
$ws.Range("AG3").Value = $false

$ws.Cells.Item(3, 46).NumberFormat = "General"  # AT3 (Bestämningsår)

$ws.Range("AW3").Value = "Liv Vikingson"
$ws.Range("AX3").Value = "Liv Vikingson"

$ws.Cells.Item(3, 51).NumberFormat = "General"  # AY3 (Projektnamn)
